# Update gh-pages output (苏州-漫展信息.xlsx) to the newly scraped data.
#
# Two events swap places / get refreshed numbers on both the "展览" sheet
# (rows 2-3) and the "全部类型" sheet (rows 3-4, shifted by one extra
# row at the top because that sheet also carries the "演出" entry).
# A handful of other rows on each of those sheets just get their
# "想去人数" (interest count) bumped.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (index 1): rows 2 and 3 swap content ---
$wsExhibit = $wb.Worksheets.Item(1)

# Row 2 -> becomes the refreshed "Good jump" event.
$wsExhibit.Cells.Item(2, 3).Value = "苏州·Good jump ACG中秋嘉年华动漫国潮文化节"
$wsExhibit.Cells.Item(2, 5).Value = "2024.09.15 10:00-09.16 17:00"
$wsExhibit.Cells.Item(2, 6).Value = 12628
$wsExhibit.Cells.Item(2, 7).Value = 49.9
$wsExhibit.Cells.Item(2, 8).Value = "https://show.bilibili.com/platform/detail.html?id=87120"
$wsExhibit.Cells.Item(2, 9).Value = "//i0.hdslb.com/bfs/openplatform/202407/yw21E7Vn1721701909995.jpeg"

# Row 3 -> becomes the "Miracle" event (prefix dropped).
$wsExhibit.Cells.Item(3, 3).Value = "苏州·Miracle☆奇迹少女 首届Live专场"
$wsExhibit.Cells.Item(3, 5).Value = "2024.09.15 13:00-09.15 15:30"
$wsExhibit.Cells.Item(3, 6).Value = 23
$wsExhibit.Cells.Item(3, 7).Value = 35
$wsExhibit.Cells.Item(3, 8).Value = "https://show.bilibili.com/platform/detail.html?id=90964"
$wsExhibit.Cells.Item(3, 9).Value = "//i2.hdslb.com/bfs/openplatform/202408/MD7PB8gD1724123134120.jpeg"

# "想去人数" (column F) bumps elsewhere on "展览".
$wsExhibit.Cells.Item(4, 6).Value = 74
$wsExhibit.Cells.Item(5, 6).Value = 41
$wsExhibit.Cells.Item(8, 6).Value = 12528
$wsExhibit.Cells.Item(10, 6).Value = 4925
$wsExhibit.Cells.Item(11, 6).Value = 4847
$wsExhibit.Cells.Item(12, 6).Value = 162
$wsExhibit.Cells.Item(16, 6).Value = 971
$wsExhibit.Cells.Item(17, 6).Value = 3
$wsExhibit.Cells.Item(19, 6).Value = 368

# --- Sheet "全部类型" (index 4): rows 3 and 4 swap content ---
$wsAll = $wb.Worksheets.Item(4)

# Row 3 -> becomes the refreshed "Good jump" event.
$wsAll.Cells.Item(3, 3).Value = "苏州·Good jump ACG中秋嘉年华动漫国潮文化节"
$wsAll.Cells.Item(3, 5).Value = "2024.09.15 10:00-09.16 17:00"
$wsAll.Cells.Item(3, 6).Value = 12628
$wsAll.Cells.Item(3, 7).Value = 49.9
$wsAll.Cells.Item(3, 8).Value = "https://show.bilibili.com/platform/detail.html?id=87120"
$wsAll.Cells.Item(3, 9).Value = "//i0.hdslb.com/bfs/openplatform/202407/yw21E7Vn1721701909995.jpeg"

# Row 4 -> becomes the "Miracle" event (prefix dropped).
$wsAll.Cells.Item(4, 3).Value = "苏州·Miracle☆奇迹少女 首届Live专场"
$wsAll.Cells.Item(4, 5).Value = "2024.09.15 13:00-09.15 15:30"
$wsAll.Cells.Item(4, 6).Value = 23
$wsAll.Cells.Item(4, 7).Value = 35
$wsAll.Cells.Item(4, 8).Value = "https://show.bilibili.com/platform/detail.html?id=90964"
$wsAll.Cells.Item(4, 9).Value = "//i2.hdslb.com/bfs/openplatform/202408/MD7PB8gD1724123134120.jpeg"

# "想去人数" (column F) bumps elsewhere on "全部类型".
$wsAll.Cells.Item(5, 6).Value = 74
$wsAll.Cells.Item(6, 6).Value = 41
$wsAll.Cells.Item(9, 6).Value = 12528
$wsAll.Cells.Item(11, 6).Value = 4925
$wsAll.Cells.Item(12, 6).Value = 4847
$wsAll.Cells.Item(13, 6).Value = 162
$wsAll.Cells.Item(17, 6).Value = 971
$wsAll.Cells.Item(18, 6).Value = 3
$wsAll.Cells.Item(20, 6).Value = 368
